# Update projection values for several "policy line" rows after getting
# rid of redundant policy lines (source data recalculation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> ordered values for columns K..S (2023 Q1 .. 2025 Q1)
$updates = @{
    9  = @(-0.4662, -0.0489, -0.0509, -0.0429, -0.295,  -0.4274, -0.0613, -0.0499, -0.3899)
    11 = @(-0.0641, -0.0694, -0.0526, -0.0471, -0.0352, -0.0202, -0.0223, -0.0112, -0.0101)
    13 = @(-0.1115, -0.0957, -0.1426, -0.1501, -0.1383, -0.038,   0.0151,  0.0197,  0.03)
    15 = @(-0.117,   0.0557,  0.1076,  0.1645,  0.4363,  0.4298,  0.4413,  0.4039,  0.2893)
    31 = @(-1.5464, -1.6509, -0.585,  -0.5631, -0.6118, -0.1892, -0.1336,  0.0243, -0.0478)
    65 = @(-0.0693,  0.0125,  0.114,   0.081,   0.0196,  0.015,   0.0236, -0.0244,  0.0356)
    67 = @( 0.04,    0.0413,  0.0429,  0.0448,  0.0464,  0.0481,  0.0495,  0.0412,  0.0365)
    69 = @(-0.0452, -0.0879, -0.1301, -0.1319, -0.0867, -0.0414,  0.0026,  0.0076,  0.0087)
    71 = @( 0.0679,  0.1303,  0.1576,  0.1184,  0.0922,  0.0924,  0.0869,  0.061,  -0.0027)
    87 = @(-0.0746,  0.1391,  0.2986,  0.1398,  0.1308,  0.1549,  0.2022,  0.1181, -0.082)
}

$columns = @("K", "L", "M", "N", "O", "P", "Q", "R", "S")

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$row"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
